$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Description" column (B) with corrected/expanded state descriptions.
# (B4's and B12's text are unchanged, so they are intentionally left untouched.)
$ws.Range("B2").Value  = "Start State; man can cross with the goat"
$ws.Range("B3").Value  = "Man can cross with nothing"
$ws.Range("B5").Value  = "Man can cross with the goat "
$ws.Range("B6").Value  = "Man can cross with the goat "
$ws.Range("B7").Value  = "Man can cross with the cabbage (following state 3)"
$ws.Range("B8").Value  = "Man can cross with the wolf (following state 4)"
$ws.Range("B9").Value  = "Man can cross with nothing"
$ws.Range("B10").Value = "Man can cross with goat"
$ws.Range("B11").Value = "Accepting State; man crossed with goat"

# Add a footnote cell with a citation, where the book title is italicized.
$noteCell = $ws.Range("H1")
$notePrefix = "Note: This is based on the DFA on page 11 of "
$noteTitle  = "Formal Language A Practical Introduction"
$noteCell.Value = $notePrefix + $noteTitle

$startPos = $notePrefix.Length + 1
$titleLen = $noteTitle.Length
$noteCell.Characters($startPos, $titleLen).Font.Italic = $true

# Update the selection to the newly added cell.
$noteCell.Select() | Out-Null

# Set the page to portrait orientation.
$ws.PageSetup.Orientation = 1
